$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Espinaca" (Femacal de La Calera) needs to
# be inserted as row 394, pushing the existing rows 394-439 down to 395-440.
$ws.Rows.Item(394).Insert()

$ws.Cells.Item(394, 1).Value = 3
$ws.Cells.Item(394, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(394, 3).Value = "Coquimbo"
$ws.Cells.Item(394, 4).Value = 44918
$ws.Cells.Item(394, 5).Value = 5
$ws.Cells.Item(394, 6).Value = 100112012
$ws.Cells.Item(394, 7).Value = "Espinaca"
$ws.Cells.Item(394, 8).Value = "Sin especificar"
$ws.Cells.Item(394, 9).Value = "Primera"
$ws.Cells.Item(394, 10).Value = 200
$ws.Cells.Item(394, 11).Value = 4500
$ws.Cells.Item(394, 12).Value = 5000
$ws.Cells.Item(394, 13).Value = 4775
$ws.Cells.Item(394, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(394, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(394, 16).Value = 1592
$ws.Cells.Item(394, 17).Value = 3
$ws.Cells.Item(394, 18).Value = "Hortaliza"
